$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 52
$ws1.Range("F3").Value = 21503
$ws1.Range("F7").Value = 30
$ws1.Range("F8").Value = 8003
$ws1.Range("F9").Value = 559
$ws1.Range("F10").Value = 46
$ws1.Range("F15").Value = 181
$ws1.Range("F20").Value = 550
$ws1.Range("F21").Value = 85
$ws1.Range("F24").Value = 88
$ws1.Range("F26").Value = 358
$ws1.Range("F27").Value = 1203
$ws1.Range("F28").Value = 60
$ws1.Range("F30").Value = 231
$ws1.Range("F31").Value = 608
$ws1.Range("F33").Value = 150
$ws1.Range("F34").Value = 5118
$ws1.Range("F37").Value = 56
$ws1.Range("F39").Value = 13227
$ws1.Range("F40").Value = 1371
$ws1.Range("F41").Value = 145
$ws1.Range("F44").Value = 323
$ws1.Range("F45").Value = 450
$ws1.Range("F47").Value = 12

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 333

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 52
$ws4.Range("F3").Value = 21503
$ws4.Range("F5").Value = 30
$ws4.Range("F6").Value = 8003
$ws4.Range("F7").Value = 559
$ws4.Range("F8").Value = 46
$ws4.Range("F13").Value = 181
$ws4.Range("F17").Value = 550
$ws4.Range("F18").Value = 85
$ws4.Range("F21").Value = 88
$ws4.Range("F23").Value = 358
$ws4.Range("F24").Value = 1203
$ws4.Range("F25").Value = 60
$ws4.Range("F27").Value = 231
$ws4.Range("F28").Value = 333
$ws4.Range("F29").Value = 608
$ws4.Range("F32").Value = 150
$ws4.Range("F34").Value = 5118
$ws4.Range("F37").Value = 56
$ws4.Range("F39").Value = 13227
$ws4.Range("F40").Value = 1371
$ws4.Range("F41").Value = 145
$ws4.Range("F44").Value = 323
$ws4.Range("F45").Value = 450
$ws4.Range("F47").Value = 12
